$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.274.65"
$ws.Range("E2").Value = "  -5.85%  "
$ws.Range("D3").Value = "1.838.24"
$ws.Range("E3").Value = "  -5.56%  "
$ws.Range("E4").Value = "  -0.52%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "331.54"
$ws.Range("E5").Value = "  -0.24%  "
$ws.Range("E6").Value = "  -0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4611"
$ws.Range("E7").Value = "  -4.93%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3865"
$ws.Range("E8").Value = "  -6.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.81"
$ws.Range("E9").Value = "  -3.84%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07847"
$ws.Range("E10").Value = "  -4.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9664"
$ws.Range("E11").Value = "  -5.24%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.85"
$ws.Range("E12").Value = "  -9.00%  "
$ws.Range("D13").Value = "1.832.12"
$ws.Range("E13").Value = "  -6.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.745"
$ws.Range("E14").Value = "  -5.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.914"
$ws.Range("E15").Value = "  -5.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06869"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "87.00"
$ws.Range("E18").Value = "  -5.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009942"
$ws.Range("E19").Value = "  -4.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.88"
$ws.Range("E20").Value = "  -5.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.003"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "28.317.26"
$ws.Range("E22").Value = "  -5.66%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.346"
$ws.Range("E23").Value = "  -5.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.08"
$ws.Range("E24").Value = "  -7.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.164"
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("D26").Value = "2.053.45"
$ws.Range("E26").Value = "  -6.43%  "
$ws.Range("E27").Value = "  -2.21%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.24"
$ws.Range("E28").Value = "  -4.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.849"
$ws.Range("E29").Value = "  -11.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.975"
$ws.Range("E30").Value = "  -6.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "116.71"
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09334"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9401"
$ws.Range("E33").Value = "  -8.06%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.288"
$ws.Range("E34").Value = "  -6.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.453"
$ws.Range("E35").Value = "  -3.04%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.327"
$ws.Range("E36").Value = "  -6.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06055"
$ws.Range("E37").Value = "  -7.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02163"
$ws.Range("E38").Value = "  -6.00%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.156"
$ws.Range("E39").Value = "  -5.93%  "
$ws.Range("E40").Value = "  -0.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5624"
$ws.Range("E41").Value = "  -6.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.564"
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.01"
$ws.Range("E43").Value = "  -6.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1785"
$ws.Range("E44").Value = "  -4.09%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.347"
$ws.Range("E45").Value = "  -8.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.244"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("E47").Value = "  -5.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5313"
$ws.Range("E48").Value = "  -5.09%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07032"
$ws.Range("E49").Value = "  -6.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.844"
$ws.Range("E50").Value = "  -7.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "113.12"
$ws.Range("E51").Value = "  -3.98%  "
